$d = $word.ActiveDocument

# Replace all occurrences of "July 03, 2022" with "July 04, 2022"
$d.Content.Find.Execute("July 03, 2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, "July 04, 2022", 2)

# Replace "September 01, 2022" with "September 02, 2022"
$d.Content.Find.Execute("September 01, 2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, "September 02, 2022", 2)
